$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-17 Wednesday" "2025-12-18 Thursday"

Replace-Text "525×3=1575" "320×2=640"
Replace-Text "637×9=5733" "804×5=4020"
Replace-Text "553×8=4424" "308×7=2156"
Replace-Text "329×3=987" "721×3=2163"
Replace-Text "193×7=1351" "526×4=2104"

Replace-Text "387×2=774" "948×9=8532"
Replace-Text "272×4=1088" "829×3=2487"
Replace-Text "483×2=966" "656×2=1312"
Replace-Text "702×8=5616" "576×2=1152"
Replace-Text "441×5=2205" "214×7=1498"

Replace-Text "794×9=7146" "299×5=1495"
Replace-Text "578×9=5202" "424×4=1696"
Replace-Text "204×4=816" "314×5=1570"
Replace-Text "803×2=1606" "662×6=3972"
Replace-Text "501×5=2505" "823×8=6584"

Replace-Text "758×6=4548" "423×9=3807"
Replace-Text "897×7=6279" "538×6=3228"
Replace-Text "253×5=1265" "748×2=1496"
Replace-Text "907×6=5442" "863×5=4315"
Replace-Text "598×8=4784" "109×6=654"

Replace-Text "740×2=1480" "592×7=4144"
Replace-Text "905×7=6335" "365×7=2555"
Replace-Text "314×4=1256" "605×6=3630"
Replace-Text "788×6=4728" "895×6=5370"
Replace-Text "511×3=1533" "162×4=648"
